# "New Hasmap Test case"
# Update the three generated-account e-mail addresses on the
# AccountCreationData sheet and make that sheet the active tab/selection
# (moving the tabSelected flag away from Credentials).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AccountCreationData")

$ws.Range("A2").Value = "newrr34r34rtest1@gmail.com"
$ws.Range("A3").Value = "newte43534st2@gmail.com"
$ws.Range("A4").Value = "new43535test3@gmail.com"

# Make AccountCreationData the active sheet/tab and move its selection.
$ws.Activate()
$ws.Range("C9").Select()
